# Auto-generated edit script applying targeted cell updates across 8 sheets
# per the Behemoth_Profits.xlsx diff (currentAveragePrice / Leve profit recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 364.1111
$ws.Range("I28").Value = 184.14285
$ws.Range("K28").Value = 184.14285
$ws.Range("M28").Value = 300.85715

$ws.Range("H51").Value = 44476.31
$ws.Range("J51").Value = 35682.668
$ws.Range("L51").Value = 35682.668
$ws.Range("N51").Value = -36650.668

$ws.Range("H98").Value = 71430860
$ws.Range("I98").Value = 76925330
$ws.Range("K98").Value = 76925330
$ws.Range("M98").Value = -76923832

$ws.Range("H100").Value = 2287.1428
$ws.Range("I100").Value = 1147.25
$ws.Range("K100").Value = 1147.25
$ws.Range("M100").Value = -606.25

$ws.Range("H122").Value = 71430860
$ws.Range("I122").Value = 76925330
$ws.Range("K122").Value = 230775990
$ws.Range("M122").Value = -230773540

$ws.Range("H135").Value = 10395.786
$ws.Range("I135").Value = 3586.3
$ws.Range("J135").Value = 27419.5
$ws.Range("K135").Value = 32276.7
$ws.Range("L135").Value = 246775.5
$ws.Range("M135").Value = -29741.7
$ws.Range("N135").Value = -251845.5

$ws.Range("H137").Value = 3492.818
$ws.Range("J137").Value = 6938.1113
$ws.Range("L137").Value = 20814.3339
$ws.Range("N137").Value = -25914.3339

$ws.Range("H138").Value = 1590091.8
$ws.Range("I138").Value = 1366
$ws.Range("J138").Value = 2225582.2
$ws.Range("K138").Value = 4098
$ws.Range("L138").Value = 6676746.600000001
$ws.Range("M138").Value = 1042
$ws.Range("N138").Value = -6687026.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 600
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -1224

$ws.Range("H97").Value = 1338.68
$ws.Range("I97").Value = 921.1053000000001
$ws.Range("K97").Value = 921.1053000000001
$ws.Range("M97").Value = -425.1053000000001

$ws.Range("H122").Value = 3937.375
$ws.Range("I122").Value = 2166.3333
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6498.999899999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4048.999899999999
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 6469.2964
$ws.Range("I132").Value = 2833.6667
$ws.Range("K132").Value = 8501.000100000001
$ws.Range("M132").Value = -5971.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 600
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -1230

$ws.Range("H99").Value = 2849.7273
$ws.Range("I99").Value = 1991.6666
$ws.Range("K99").Value = 1991.6666
$ws.Range("M99").Value = -493.6666

$ws.Range("H134").Value = 27377.238
$ws.Range("I134").Value = 3298.2104
$ws.Range("K134").Value = 9894.6312
$ws.Range("M134").Value = -7359.6312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 75000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 75000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 75000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -76022

$ws.Range("H68").Value = 61242.25
$ws.Range("J68").Value = 56663
$ws.Range("L68").Value = 56663
$ws.Range("N68").Value = -58161

$ws.Range("H71").Value = 61242.25
$ws.Range("J71").Value = 56663
$ws.Range("L71").Value = 169989
$ws.Range("N71").Value = -177477

$ws.Range("H122").Value = 3984.889
$ws.Range("I122").Value = 2699.8462
$ws.Range("K122").Value = 8099.5386
$ws.Range("M122").Value = -5649.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 8771
$ws.Range("I103").Value = 1899.5
$ws.Range("K103").Value = 5698.5
$ws.Range("M103").Value = -4819.5

$ws.Range("H107").Value = 437.8
$ws.Range("I107").Value = 355
$ws.Range("J107").Value = 467.9091
$ws.Range("K107").Value = 1065
$ws.Range("L107").Value = 1403.7273
$ws.Range("M107").Value = 855
$ws.Range("N107").Value = -5243.7273

$ws.Range("H109").Value = 2976.6667
$ws.Range("I109").Value = 2915
$ws.Range("K109").Value = 8745
$ws.Range("M109").Value = -7705

$ws.Range("H132").Value = 2225.1765
$ws.Range("J132").Value = 3118.5
$ws.Range("L132").Value = 28066.5
$ws.Range("N132").Value = -33126.5

$ws.Range("H138").Value = 1099.3334
$ws.Range("I138").Value = 1099.3334
$ws.Range("K138").Value = 3298.0002
$ws.Range("M138").Value = 1841.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.588234
$ws.Range("I2").Value = 80.75
$ws.Range("J2").Value = 78.55556
$ws.Range("K2").Value = 80.75
$ws.Range("L2").Value = 78.55556
$ws.Range("M2").Value = 32.25
$ws.Range("N2").Value = -304.55556

$ws.Range("H102").Value = 3070.76
$ws.Range("I102").Value = 2615.1333
$ws.Range("J102").Value = 3754.2
$ws.Range("K102").Value = 2615.1333
$ws.Range("L102").Value = 3754.2
$ws.Range("M102").Value = -993.1333
$ws.Range("N102").Value = -6998.2

$ws.Range("H105").Value = 85607.836
$ws.Range("J105").Value = 85607.836
$ws.Range("L105").Value = 85607.836
$ws.Range("N105").Value = -92595.836

$ws.Range("H126").Value = 4760.6924
$ws.Range("I126").Value = 4999.4
$ws.Range("J126").Value = 4611.5
$ws.Range("K126").Value = 14998.2
$ws.Range("L126").Value = 13834.5
$ws.Range("M126").Value = -12528.2
$ws.Range("N126").Value = -18774.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3905.4167
$ws.Range("I100").Value = 4385
$ws.Range("K100").Value = 4385
$ws.Range("M100").Value = -3844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 121000
$ws.Range("J57").Value = 121000
$ws.Range("L57").Value = 121000
$ws.Range("N57").Value = -122508

$ws.Range("H81").Value = 67267.336
$ws.Range("I81").Value = 100400
$ws.Range("K81").Value = 200800
$ws.Range("M81").Value = -199739

$ws.Range("H84").Value = 67267.336
$ws.Range("I84").Value = 100400
$ws.Range("K84").Value = 1004000
$ws.Range("M84").Value = -998696

$ws.Range("H122").Value = 4599.05
$ws.Range("I122").Value = 3225.4666
$ws.Range("J122").Value = 8719.799999999999
$ws.Range("K122").Value = 9676.399800000001
$ws.Range("L122").Value = 26159.4
$ws.Range("M122").Value = -7226.399800000001
$ws.Range("N122").Value = -31059.4

$ws.Range("H127").Value = 62181.668
$ws.Range("J127").Value = 64618
$ws.Range("L127").Value = 64618
$ws.Range("N127").Value = -74538

$ws.Range("H141").Value = 58996.5
$ws.Range("J141").Value = 58996.5
$ws.Range("L141").Value = 58996.5
$ws.Range("N141").Value = -69356.5
